$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.291.47"
$ws.Range("E2").Value = "  -3.33%  "
$ws.Range("D3").Value = "2.466.99"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'312.71"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").Value = "'94.29"
$ws.Range("E6").Value = "  -6.67%  "
$ws.Range("D7").Value = "'0.552"
$ws.Range("E7").Value = "  -3.14%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.502"
$ws.Range("E9").Value = "  -4.45%  "
$ws.Range("D10").Value = "'33.57"
$ws.Range("E10").Value = "  -6.56%  "
$ws.Range("D11").Value = "'0.0780"
$ws.Range("E11").Value = "  -2.89%  "
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "'7.03"
$ws.Range("E13").Value = "  -3.55%  "
$ws.Range("D14").Value = "2.848.13"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "2.464.86"
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("E16").Value = "  -6.22%  "
$ws.Range("D17").Value = "'0.786"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").Value = "41.261.47"
$ws.Range("E18").Value = "  -3.31%  "
$ws.Range("D19").Value = "'6.34"
$ws.Range("E19").Value = "  -5.88%  "
$ws.Range("D20").Value = "0.0₃0921"
$ws.Range("E20").Value = "  -3.03%  "
$ws.Range("D21").Value = "'11.50"
$ws.Range("E21").Value = "  -5.39%  "
$ws.Range("D22").Value = "'68.05"
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("D23").Value = "'237.35"
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("D25").Value = "'1.93"
$ws.Range("E25").Value = "  -5.17%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'24.46"
$ws.Range("E27").Value = "  -6.60%  "
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("D29").Value = "'9.71"
$ws.Range("E29").Value = "  -4.23%  "
$ws.Range("D30").Value = "'36.03"
$ws.Range("E30").Value = "  -7.86%  "
$ws.Range("D31").Value = "'153.49"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("D32").Value = "'5.57"
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("D33").Value = "'2.61"
$ws.Range("E33").Value = "  -6.15%  "
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Value = "'0.0756"
$ws.Range("E35").Value = "  -3.80%  "
$ws.Range("D36").Value = "'3.01"
$ws.Range("E36").Value = "  -5.49%  "
$ws.Range("D37").Value = "'1.89"
$ws.Range("E37").Value = "  -6.80%  "
$ws.Range("D38").Value = "'16.95"
$ws.Range("E38").Value = "  -6.94%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.104"
$ws.Range("E39").Value = "  -6.70%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.114"
$ws.Range("E40").Value = "  -3.78%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'4.32"
$ws.Range("E41").Value = "  +2.32%  "
$ws.Range("D42").Value = "'21.28"
$ws.Range("E42").Value = "  -2.70%  "
$ws.Range("D43").Value = "'1.01"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "1.985.02"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "'0.0286"
$ws.Range("E45").Value = "  -4.55%  "
$ws.Range("D46").Value = "'3.07"
$ws.Range("E46").Value = "  -6.37%  "
$ws.Range("D47").Value = "'8.72"
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").Value = "'69.90"
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("D49").Value = "'76.16"
$ws.Range("E49").Value = "  -5.11%  "
$ws.Range("D50").Value = "'97.10"
$ws.Range("E50").Value = "  -3.81%  "
$ws.Range("D51").Value = "'0.179"
$ws.Range("E51").Value = "  -5.78%  "
